$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16502911924804595"
$ws1.Range("B2").Value = "go_stims-1650291192428458.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291192447459.csv"
$ws1.Range("B4").Value = "go_stims-16502911924494581.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911924794586.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502911941134646"
$ws2.Range("B2").Value = "TB-16502911940894568.csv"
$ws2.Range("B3").Value = "ZB-match_3-1650291192624456.csv"
$ws2.Range("B4").Value = "TB-16502911933174572.csv"
$ws2.Range("B5").Value = "OB-16502911931834574.csv"
$ws2.Range("B6").Value = "OB-16502911927814565.csv"
$ws2.Range("B7").Value = "ZB-match_4-16502911925754573.csv"
$ws2.Range("B8").Value = "OB-16502911927594569.csv"
$ws2.Range("B9").Value = "TB-16502911938944657.csv"
$ws2.Range("B10").Value = "ZB-match_9-16502911925444615.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650291194116457"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1650291194194458"
$ws4.Range("B2").Value = "MM_stims-1650291194144462.csv"
$ws4.Range("B3").Value = "ZM_stims-1650291194118466.csv"
$ws4.Range("B4").Value = "MM_stims-16502911941754632.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911941464565.csv"
$ws4.Range("B6").Value = "MM_stims-1650291194191458.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911941774597.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650291194270457"
$ws5.Range("B2").Value = "SAT_stims-16502911941994593.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911942564585.csv"
$ws5.Range("B4").Value = "SAT_stims-16502911942244577.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911942394576.csv"
